{"js": "// Apply the Annahmen.docx text edits described by the diff.\n// Only genuine visible-text changes are applied; the diff's run-splitting\n// and <w:proofErr> additions are spell/grammar-check markup artifacts\n// that do not alter the document's visible text, so they are skipped.\n\nconst body = context.document.body;\n\n// 1) \"...m\u00fcssen exakt mit den Spezifizierungen die zu erkennen sind\n//    \u00fcbereinstimmen.\" -> add two commas.\nconst r1 = body.search(\n  \"Titel und Anreden m\u00fcssen exakt mit den Spezifizierungen die zu erkennen sind \u00fcbereinstimmen.\",\n  { matchCase: true }\n);\nr1.load(\"items\");\nawait context.sync();\nif (r1.items.length > 0) {\n  r1.items[0].insertText(\n    \"Titel und Anreden m\u00fcssen exakt mit den Spezifizierungen, die zu erkennen sind, \u00fcbereinstimmen.\",\n    \"Replace\"\n  );\n}\nawait context.sync();\n\n// 2) Remove \" und ausschlie\u00dflich Buchstaben des Alphabets (A-Z)\" so the\n//    paragraph ends after \"eingegeben\".\nconst r2 = body.search(\n  \"Gro\u00df- und Kleinschreibung korrekt eingegeben und ausschlie\u00dflich Buchstaben des Alphabets (A-Z)\",\n  { matchCase: true }\n);\nr2.load(\"items\");\nawait context.sync();\nif (r2.items.length > 0) {\n  r2.items[0].insertText(\"Gro\u00df- und Kleinschreibung korrekt eingegeben\", \"Replace\");\n}\nawait context.sync();\n\n// 3) Append a new example sentence at the end of the \"Adelstitel\" bullet.\nconst r3 = body.search(\n  \"Adelstitel werden in den Vornamen geschrieben (M\u00f6glichkeit zur Bearbeitung des Vor- und Nachnamens kann danach genutzt werden)\",\n  { matchCase: true }\n);\nr3.load(\"items\");\nawait context.sync();\nif (r3.items.length > 0) {\n  r3.items[0].insertText(\" (Beispiel. Max Freiherr von Waldenbuch)\", \"End\");\n}\nawait context.sync();\n\n// 4) \"...hinzugef\u00fcgt werden\" -> \"...hinzugef\u00fcgt/bearbeitet werden.\" on the\n//    \"Vorschau\" (Anrede/Titel) bullet.\nconst r4 = body.search(\n  \"Diese m\u00fcssen vor Abspeicherung des Kontakts \u00fcber die entsprechenden Fenster hinzugef\u00fcgt werden\",\n  { matchCase: true }\n);\nr4.load(\"items\");\nawait context.sync();\nif (r4.items.length > 0) {\n  r4.items[0].insertText(\n    \"Diese m\u00fcssen vor Abspeicherung des Kontakts \u00fcber die entsprechenden Fenster hinzugef\u00fcgt/bearbeitet werden.\",\n    \"Replace\"\n  );\n}\nawait context.sync();\n", "ps1": "# Apply the Annahmen.docx text edits described by the diff.\n# Only genuine visible-text changes are applied; the diff's run-splitting\n# and proofErr (spell/grammar-check) additions are markup artifacts that\n# do not alter the document's visible text, so they are skipped here.\n\n$d = $word.ActiveDocument\n$wdFindContinue = 1\n$wdReplaceAll = 2\n$wdCollapseEnd = 0\n\n# 1) \"...m\u00fcssen exakt mit den Spezifizierungen die zu erkennen sind\n#    \u00fcbereinstimmen.\" -> add two commas.\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"Titel und Anreden m\u00fcssen exakt mit den Spezifizierungen die zu erkennen sind \u00fcbereinstimmen.\"\n$find.Replacement.Text = \"Titel und Anreden m\u00fcssen exakt mit den Spezifizierungen, die zu erkennen sind, \u00fcbereinstimmen.\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $find.Replacement.Text, $wdReplaceAll)\n\n# 2) Remove \" und ausschlie\u00dflich Buchstaben des Alphabets (A-Z)\" so the\n#    paragraph ends after \"eingegeben\".\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"Gro\u00df- und Kleinschreibung korrekt eingegeben und ausschlie\u00dflich Buchstaben des Alphabets (A-Z)\"\n$find.Replacement.Text = \"Gro\u00df- und Kleinschreibung korrekt eingegeben\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $find.Replacement.Text, $wdReplaceAll)\n\n# 3) Append a new example sentence at the end of the \"Adelstitel\" bullet.\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"Adelstitel werden in den Vornamen geschrieben (M\u00f6glichkeit zur Bearbeitung des Vor- und Nachnamens kann danach genutzt werden)\"\n$found = $find.Execute()\nif ($found) {\n    $rng = $find.Parent\n    $rng.Collapse($wdCollapseEnd)\n    $rng.InsertAfter(\" (Beispiel. Max Freiherr von Waldenbuch)\")\n}\n\n# 4) \"...hinzugef\u00fcgt werden\" -> \"...hinzugef\u00fcgt/bearbeitet werden.\" on the\n#    \"Vorschau\" (Anrede/Titel) bullet.\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"Diese m\u00fcssen vor Abspeicherung des Kontakts \u00fcber die entsprechenden Fenster hinzugef\u00fcgt werden\"\n$find.Replacement.Text = \"Diese m\u00fcssen vor Abspeicherung des Kontakts \u00fcber die entsprechenden Fenster hinzugef\u00fcgt/bearbeitet werden.\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $find.Replacement.Text, $wdReplaceAll)\n"}
